$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row to the new machine-friendly column names
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "mx_state"
$ws.Range("B1").Value2 = "mx_municipality"
$ws.Range("C1").Value2 = "n_matriculas"
$ws.Range("D1").Value2 = "pct_matriculas"

# ---------------------------------------------------------------------------
# 2) Title-case the little connector words ("de", "del", "el", "los", "las",
#    "la", "y") inside the state/municipality names in columns A and B, plus
#    a couple of one-off spelling fixes ("GUANAJUATO" -> "Guanajuato" and
#    "MonteMorelos" -> "Montemorelos").
# ---------------------------------------------------------------------------
$words = @("de", "del", "el", "los", "las", "la", "y")

function Transform-Text($s) {
    if ($s.Equals("GUANAJUATO")) { return "Guanajuato" }
    if ($s.Equals("MonteMorelos")) { return "Montemorelos" }

    $parts = $s.Split(" ")
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($words -ccontains $parts[$i]) {
            $parts[$i] = $parts[$i].Substring(0, 1).ToUpper() + $parts[$i].Substring(1)
        }
    }
    return [string]::Join(" ", $parts)
}

$lastRow = 1094
for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in @(1, 2)) {
        $cell = $ws.Cells.Item($r, $col)
        $v = $cell.Value2
        if ($v -ne $null -and $v.GetType().Name -eq "String") {
            $nv = Transform-Text $v
            if (-not $v.Equals($nv)) {
                $cell.Value2 = $nv
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Tiny floating point re-roundings in the "pct_matriculas" column that the
#    refreshed pipeline produced (one ULP differences from 7/7206 and
#    652/7206). Done by direct cell assignment (not Range.Replace, which
#    round-trips the number through text and can perturb the last bit).
# ---------------------------------------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $v = $cell.Value2
    if ($v -ne $null) {
        if ($v.Equals(0.0009714127116291979)) {
            $cell.Value2 = 0.000971412711629198
        } elseif ($v.Equals(0.09048015542603387)) {
            $cell.Value2 = 0.09048015542603388
        }
    }
}

# ---------------------------------------------------------------------------
# 4) Drop the trailing footnote rows (sample size / source / author / date)
#    that used to live below the grand-total row.
# ---------------------------------------------------------------------------
$ws.Rows("1096:1100").Delete() | Out-Null
